$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the scientific (Latin binomial) portion of each species' common name,
# which lives inside the parentheses in column A, with markdown-style
# asterisks (italics markers) -- e.g. "(Felis nigripes)" -> "(*Felis nigripes*)".
# All other columns / values are left untouched.

$ws.Range("A2").Value = "Black-footed cat (*Felis nigripes*)"
$ws.Range("A3").Value = "Fynbos golden mole (*Amblysomus corriae*)"
$ws.Range("A4").Value = "Namaqua dune mole-rat (*Bathyergus janetta*)"
$ws.Range("A5").Value = "Riverine rabbit (*Bunolagus monticularis*)"
$ws.Range("A6").Value = "Bontebok (*Damaliscus pygargus pygargus*)"
$ws.Range("A7").Value = "Cape mountain zebra (*Equus zebra zebra*)"
$ws.Range("A8").Value = "Bush-tailed hairy-footed gerbil (*Gerbillurus vallinus*)"
$ws.Range("A9").Value = "Spectacled dormouse (*Graphiurus acularis*)"
$ws.Range("A10").Value = "Verreaux's mouse (*Myomyscus verreauxii*)"
$ws.Range("A11").Value = "Vlei rat (*Otomys auratus*)"
$ws.Range("A12").Value = "Fynbos vlei rat (*Otomys irroratus*)"
$ws.Range("A13").Value = "Springhare (*Pedetes capensis*)"
$ws.Range("A14").Value = "Barbour's rock mouse (*Petromyscus barbouri*)"
$ws.Range("A15").Value = "Grysbok (*Raphicerus melanotis*)"
